$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": append two new rows (File Name = Path-And-Name pairs
# for the two newly-generated handoff files), growing the table from
# A1:G3 to A1:G5.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$rowO1 = $loOverview.ListRows.Add()
$r = $rowO1.Range.Row
$wsOverview.Cells.Item($r, 1).Value = "aa56937a-9016-49a7-a98f-a0b9fa787b97.md"
$wsOverview.Cells.Item($r, 2).Value = "e2e\aa56937a-9016-49a7-a98f-a0b9fa787b97.md"
$wsOverview.Cells.Item($r, 3).Value = ".md"
$wsOverview.Cells.Item($r, 4).Value = "'"
$wsOverview.Cells.Item($r, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 7).Value = "2016-09-03 00:44:26"
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($r, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa56937a9016/e2e/aa56937a-9016-49a7-a98f-a0b9fa787b97.md", "", "", "e2e\aa56937a-9016-49a7-a98f-a0b9fa787b97.md")

$rowO2 = $loOverview.ListRows.Add()
$r = $rowO2.Range.Row
$wsOverview.Cells.Item($r, 1).Value = "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md"
$wsOverview.Cells.Item($r, 2).Value = "e2e\f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md"
$wsOverview.Cells.Item($r, 3).Value = ".md"
$wsOverview.Cells.Item($r, 4).Value = "'"
$wsOverview.Cells.Item($r, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 7).Value = "2016-09-03 00:44:26"
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($r, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a282bf3bc4/e2e/f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md", "", "", "e2e\f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn": append the same two handoff files as new table rows,
# growing the table from A1:P3 to A1:P5.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$rowZ1 = $loZh.ListRows.Add()
$r = $rowZ1.Range.Row
$wsZh.Cells.Item($r, 1).Value = "aa56937a-9016-49a7-a98f-a0b9fa787b97.md"
$wsZh.Cells.Item($r, 2).Value = ".md"
$wsZh.Cells.Item($r, 3).Value = "Ready for handoff"
$wsZh.Cells.Item($r, 4).Value = "e2e"
$wsZh.Cells.Item($r, 5).Value = "ht"
$wsZh.Cells.Item($r, 6).Value = "'False"
$wsZh.Cells.Item($r, 7).Value = "aa56937a-9016-49a7-a98f-a0b9fa787b97.a71d78b3800ab8a06d2b09692555633323154683.zh-cn.xlf"
$wsZh.Cells.Item($r, 8).Value = "2016-09-03 00:44:21"
$wsZh.Cells.Item($r, 9).Value = "'"
$wsZh.Cells.Item($r, 10).Value = "'"
$wsZh.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($r, 12).Value = "'"
$wsZh.Cells.Item($r, 13).Value = "'True"
$wsZh.Cells.Item($r, 14).Value = "'"
$wsZh.Cells.Item($r, 15).Value = "'False"
$wsZh.Cells.Item($r, 16).Value = "'"
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa56937a9016/e2e/aa56937a-9016-49a7-a98f-a0b9fa787b97.md", "", "", "aa56937a-9016-49a7-a98f-a0b9fa787b97.md")

$rowZ2 = $loZh.ListRows.Add()
$r = $rowZ2.Range.Row
$wsZh.Cells.Item($r, 1).Value = "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md"
$wsZh.Cells.Item($r, 2).Value = ".md"
$wsZh.Cells.Item($r, 3).Value = "Ready for handoff"
$wsZh.Cells.Item($r, 4).Value = "e2e"
$wsZh.Cells.Item($r, 5).Value = "ht"
$wsZh.Cells.Item($r, 6).Value = "'False"
$wsZh.Cells.Item($r, 7).Value = "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.8da6902b86b95236c1115537909eb6e291a5ac00.zh-cn.xlf"
$wsZh.Cells.Item($r, 8).Value = "2016-09-03 00:44:21"
$wsZh.Cells.Item($r, 9).Value = "'"
$wsZh.Cells.Item($r, 10).Value = "'"
$wsZh.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($r, 12).Value = "'"
$wsZh.Cells.Item($r, 13).Value = "'True"
$wsZh.Cells.Item($r, 14).Value = "'"
$wsZh.Cells.Item($r, 15).Value = "'False"
$wsZh.Cells.Item($r, 16).Value = "'"
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a282bf3bc4/e2e/f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md", "", "", "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md")

# ---------------------------------------------------------------------
# Sheet "de-de": append the same two handoff files as new table rows,
# growing the table from A1:P3 to A1:P5.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$rowD1 = $loDe.ListRows.Add()
$r = $rowD1.Range.Row
$wsDe.Cells.Item($r, 1).Value = "aa56937a-9016-49a7-a98f-a0b9fa787b97.md"
$wsDe.Cells.Item($r, 2).Value = ".md"
$wsDe.Cells.Item($r, 3).Value = "Ready for handoff"
$wsDe.Cells.Item($r, 4).Value = "e2e"
$wsDe.Cells.Item($r, 5).Value = "ht"
$wsDe.Cells.Item($r, 6).Value = "'False"
$wsDe.Cells.Item($r, 7).Value = "aa56937a-9016-49a7-a98f-a0b9fa787b97.a71d78b3800ab8a06d2b09692555633323154683.de-de.xlf"
$wsDe.Cells.Item($r, 8).Value = "2016-09-03 00:44:26"
$wsDe.Cells.Item($r, 9).Value = "'"
$wsDe.Cells.Item($r, 10).Value = "'"
$wsDe.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($r, 12).Value = "'"
$wsDe.Cells.Item($r, 13).Value = "'True"
$wsDe.Cells.Item($r, 14).Value = "'"
$wsDe.Cells.Item($r, 15).Value = "'False"
$wsDe.Cells.Item($r, 16).Value = "'"
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa56937a9016/e2e/aa56937a-9016-49a7-a98f-a0b9fa787b97.md", "", "", "aa56937a-9016-49a7-a98f-a0b9fa787b97.md")

$rowD2 = $loDe.ListRows.Add()
$r = $rowD2.Range.Row
$wsDe.Cells.Item($r, 1).Value = "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md"
$wsDe.Cells.Item($r, 2).Value = ".md"
$wsDe.Cells.Item($r, 3).Value = "Ready for handoff"
$wsDe.Cells.Item($r, 4).Value = "e2e"
$wsDe.Cells.Item($r, 5).Value = "ht"
$wsDe.Cells.Item($r, 6).Value = "'False"
$wsDe.Cells.Item($r, 7).Value = "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.8da6902b86b95236c1115537909eb6e291a5ac00.de-de.xlf"
$wsDe.Cells.Item($r, 8).Value = "2016-09-03 00:44:26"
$wsDe.Cells.Item($r, 9).Value = "'"
$wsDe.Cells.Item($r, 10).Value = "'"
$wsDe.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($r, 12).Value = "'"
$wsDe.Cells.Item($r, 13).Value = "'True"
$wsDe.Cells.Item($r, 14).Value = "'"
$wsDe.Cells.Item($r, 15).Value = "'False"
$wsDe.Cells.Item($r, 16).Value = "'"
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a282bf3bc4/e2e/f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md", "", "", "f2a282bf-3bc4-4f64-8c82-97f61d2fea96.md")
